# Backup QR Scanner data - shifts log rows 21-111 up by one entry
# (new scans were recorded, pushing the window forward) and drops the
# now-stale last row (112). Columns B (Subject), C (Log Date) and F
# (User) are constant across the whole sheet, so only A (Student ID),
# D (Log Time) and E (Type) need to be rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (Student ID) to remain text for numeric-looking values
# (these IDs must stay text, matching the sheet's existing
# numberStoredAsText convention, instead of being auto-coerced to Number).
$ws.Range("A21:A111").NumberFormat = "@"

$ws.Range("A21").Value = "201563"
$ws.Range("D21").Value = "11:28:13"
$ws.Range("E21").Value = "Manual"
$ws.Range("A22").Value = "180804"
$ws.Range("D22").Value = "11:28:23"
$ws.Range("E22").Value = "Manual"
$ws.Range("A23").Value = "211167"
$ws.Range("D23").Value = "11:28:31"
$ws.Range("E23").Value = "Manual"
$ws.Range("A24").Value = "211169"
$ws.Range("D24").Value = "11:28:41"
$ws.Range("E24").Value = "Manual"
$ws.Range("A25").Value = "211079"
$ws.Range("D25").Value = "11:28:48"
$ws.Range("E25").Value = "Manual"
$ws.Range("A26").Value = "200468"
$ws.Range("D26").Value = "11:29:01"
$ws.Range("E26").Value = "Manual"
$ws.Range("A27").Value = "191109"
$ws.Range("D27").Value = "11:29:08"
$ws.Range("E27").Value = "Manual"
$ws.Range("A28").Value = "191088"
$ws.Range("D28").Value = "11:29:17"
$ws.Range("E28").Value = "Manual"
$ws.Range("A29").Value = "200423"
$ws.Range("D29").Value = "11:29:25"
$ws.Range("E29").Value = "Manual"
$ws.Range("A30").Value = "190931"
$ws.Range("D30").Value = "11:29:35"
$ws.Range("E30").Value = "Manual"
$ws.Range("A31").Value = "190922"
$ws.Range("D31").Value = "11:29:50"
$ws.Range("E31").Value = "Manual"
$ws.Range("A32").Value = "191131"
$ws.Range("D32").Value = "11:30:02"
$ws.Range("E32").Value = "Manual"
$ws.Range("A33").Value = "190803"
$ws.Range("D33").Value = "11:30:15"
$ws.Range("E33").Value = "Manual"
$ws.Range("A34").Value = "181013"
$ws.Range("D34").Value = "11:30:22"
$ws.Range("E34").Value = "Manual"
$ws.Range("A35").Value = "201495"
$ws.Range("D35").Value = "11:30:50"
$ws.Range("E35").Value = "Manual"
$ws.Range("A36").Value = "201026"
$ws.Range("D36").Value = "11:33:15"
$ws.Range("E36").Value = "Manual"
$ws.Range("A37").Value = "200850"
$ws.Range("D37").Value = "11:33:26"
$ws.Range("E37").Value = "Manual"
$ws.Range("A38").Value = "200866"
$ws.Range("D38").Value = "11:33:36"
$ws.Range("E38").Value = "Manual"
$ws.Range("A39").Value = "211118"
$ws.Range("D39").Value = "11:36:25"
$ws.Range("E39").Value = "Manual"
$ws.Range("A40").Value = "200904"
$ws.Range("D40").Value = "11:37:01"
$ws.Range("E40").Value = "Manual"
$ws.Range("A41").Value = "201632"
$ws.Range("D41").Value = "11:37:08"
$ws.Range("E41").Value = "Manual"
$ws.Range("A42").Value = "201190"
$ws.Range("D42").Value = "11:37:15"
$ws.Range("E42").Value = "Manual"
$ws.Range("A43").Value = "200824"
$ws.Range("D43").Value = "11:37:21"
$ws.Range("E43").Value = "Manual"
$ws.Range("A44").Value = "201197"
$ws.Range("D44").Value = "11:37:28"
$ws.Range("E44").Value = "Manual"
$ws.Range("A45").Value = "200914"
$ws.Range("D45").Value = "11:37:39"
$ws.Range("E45").Value = "Manual"
$ws.Range("A46").Value = "201065"
$ws.Range("D46").Value = "11:38:11"
$ws.Range("E46").Value = "Manual"
$ws.Range("A47").Value = "191478"
$ws.Range("D47").Value = "11:38:32"
$ws.Range("E47").Value = "Manual"
$ws.Range("A48").Value = "200999"
$ws.Range("D48").Value = "11:38:40"
$ws.Range("E48").Value = "Manual"
$ws.Range("A49").Value = "201157"
$ws.Range("D49").Value = "11:38:51"
$ws.Range("E49").Value = "Manual"
$ws.Range("A50").Value = "190314"
$ws.Range("D50").Value = "11:39:02"
$ws.Range("E50").Value = "Manual"
$ws.Range("A51").Value = "202162"
$ws.Range("D51").Value = "11:39:12"
$ws.Range("E51").Value = "Manual"
$ws.Range("A52").Value = "201819"
$ws.Range("D52").Value = "11:39:21"
$ws.Range("E52").Value = "Manual"
$ws.Range("A53").Value = "201990"
$ws.Range("D53").Value = "11:39:29"
$ws.Range("E53").Value = "Manual"
$ws.Range("A54").Value = "211175"
$ws.Range("D54").Value = "11:39:42"
$ws.Range("E54").Value = "Manual"
$ws.Range("A55").Value = "201795"
$ws.Range("D55").Value = "11:39:58"
$ws.Range("E55").Value = "Manual"
$ws.Range("A56").Value = "211174"
$ws.Range("D56").Value = "11:41:21"
$ws.Range("E56").Value = "Manual"
$ws.Range("A57").Value = "211092"
$ws.Range("D57").Value = "11:41:29"
$ws.Range("E57").Value = "Manual"
$ws.Range("A58").Value = "211046"
$ws.Range("D58").Value = "11:41:38"
$ws.Range("E58").Value = "Manual"
$ws.Range("A59").Value = "211242"
$ws.Range("D59").Value = "11:41:50"
$ws.Range("E59").Value = "Manual"
$ws.Range("A60").Value = "211010"
$ws.Range("D60").Value = "11:41:58"
$ws.Range("E60").Value = "Manual"
$ws.Range("A61").Value = "190968"
$ws.Range("D61").Value = "11:42:09"
$ws.Range("E61").Value = "Manual"
$ws.Range("A62").Value = "200933"
$ws.Range("D62").Value = "11:42:20"
$ws.Range("E62").Value = "Manual"
$ws.Range("A63").Value = "201825"
$ws.Range("D63").Value = "11:42:35"
$ws.Range("E63").Value = "Manual"
$ws.Range("A64").Value = "190801"
$ws.Range("D64").Value = "11:42:42"
$ws.Range("E64").Value = "Manual"
$ws.Range("A65").Value = "201465"
$ws.Range("D65").Value = "11:42:49"
$ws.Range("E65").Value = "Scan"
$ws.Range("A66").Value = "201171"
$ws.Range("D66").Value = "11:42:52"
$ws.Range("E66").Value = "Scan"
$ws.Range("A67").Value = "200491"
$ws.Range("D67").Value = "11:42:54"
$ws.Range("E67").Value = "Scan"
$ws.Range("A68").Value = "200490"
$ws.Range("D68").Value = "11:42:56"
$ws.Range("E68").Value = "Scan"
$ws.Range("A69").Value = "200228"
$ws.Range("D69").Value = "11:42:59"
$ws.Range("E69").Value = "Scan"
$ws.Range("A70").Value = "201669"
$ws.Range("D70").Value = "11:43:01"
$ws.Range("E70").Value = "Scan"
$ws.Range("A71").Value = "202004"
$ws.Range("D71").Value = "11:43:03"
$ws.Range("E71").Value = "Scan"
$ws.Range("A72").Value = "200116"
$ws.Range("D72").Value = "11:43:05"
$ws.Range("E72").Value = "Scan"
$ws.Range("A73").Value = "190796"
$ws.Range("D73").Value = "11:43:07"
$ws.Range("E73").Value = "Scan"
$ws.Range("A74").Value = "190981"
$ws.Range("D74").Value = "11:43:09"
$ws.Range("E74").Value = "Scan"
$ws.Range("A75").Value = "200163"
$ws.Range("D75").Value = "11:43:12"
$ws.Range("E75").Value = "Scan"
$ws.Range("A76").Value = "191061"
$ws.Range("D76").Value = "11:43:15"
$ws.Range("E76").Value = "Scan"
$ws.Range("A77").Value = "200744"
$ws.Range("D77").Value = "11:43:17"
$ws.Range("E77").Value = "Scan"
$ws.Range("A78").Value = "201564"
$ws.Range("D78").Value = "11:43:20"
$ws.Range("E78").Value = "Scan"
$ws.Range("A79").Value = "200804"
$ws.Range("D79").Value = "11:43:22"
$ws.Range("E79").Value = "Scan"
$ws.Range("A80").Value = "200792"
$ws.Range("D80").Value = "11:43:24"
$ws.Range("E80").Value = "Scan"
$ws.Range("A81").Value = "200628"
$ws.Range("D81").Value = "11:43:26"
$ws.Range("E81").Value = "Scan"
$ws.Range("A82").Value = "211216"
$ws.Range("D82").Value = "11:43:28"
$ws.Range("E82").Value = "Scan"
$ws.Range("A83").Value = "211197"
$ws.Range("D83").Value = "11:43:31"
$ws.Range("E83").Value = "Scan"
$ws.Range("A84").Value = "201397"
$ws.Range("D84").Value = "11:43:33"
$ws.Range("E84").Value = "Scan"
$ws.Range("A85").Value = "200917"
$ws.Range("D85").Value = "11:43:35"
$ws.Range("E85").Value = "Scan"
$ws.Range("A86").Value = "201051"
$ws.Range("D86").Value = "11:43:37"
$ws.Range("E86").Value = "Scan"
$ws.Range("A87").Value = "201501"
$ws.Range("D87").Value = "11:43:40"
$ws.Range("E87").Value = "Scan"
$ws.Range("A88").Value = "200905"
$ws.Range("D88").Value = "11:43:42"
$ws.Range("E88").Value = "Scan"
$ws.Range("A89").Value = "200938"
$ws.Range("D89").Value = "11:43:45"
$ws.Range("E89").Value = "Scan"
$ws.Range("A90").Value = "211043"
$ws.Range("D90").Value = "11:43:47"
$ws.Range("E90").Value = "Scan"
$ws.Range("A91").Value = "211133"
$ws.Range("D91").Value = "11:43:49"
$ws.Range("E91").Value = "Scan"
$ws.Range("A92").Value = "211741"
$ws.Range("D92").Value = "11:43:50"
$ws.Range("E92").Value = "Scan"
$ws.Range("A93").Value = "211245"
$ws.Range("D93").Value = "11:43:52"
$ws.Range("E93").Value = "Scan"
$ws.Range("A94").Value = "190929"
$ws.Range("D94").Value = "11:43:55"
$ws.Range("E94").Value = "Scan"
$ws.Range("A95").Value = "191119"
$ws.Range("D95").Value = "11:43:56"
$ws.Range("E95").Value = "Scan"
$ws.Range("A96").Value = "211210"
$ws.Range("D96").Value = "11:44:43"
$ws.Range("E96").Value = "Manual"
$ws.Range("A97").Value = "201880"
$ws.Range("D97").Value = "11:45:28"
$ws.Range("E97").Value = "Manual"
$ws.Range("A98").Value = "211246"
$ws.Range("D98").Value = "11:45:59"
$ws.Range("E98").Value = "Manual"
$ws.Range("A99").Value = "201823"
$ws.Range("D99").Value = "11:46:33"
$ws.Range("E99").Value = "Manual"
$ws.Range("A100").Value = "200928"
$ws.Range("D100").Value = "11:47:11"
$ws.Range("E100").Value = "Manual"
$ws.Range("A101").Value = "201237"
$ws.Range("D101").Value = "11:47:35"
$ws.Range("E101").Value = "Manual"
$ws.Range("A102").Value = "201513"
$ws.Range("D102").Value = "11:48:02"
$ws.Range("E102").Value = "Manual"
$ws.Range("A103").Value = "200997"
$ws.Range("D103").Value = "11:48:26"
$ws.Range("E103").Value = "Manual"
$ws.Range("A104").Value = "181004"
$ws.Range("D104").Value = "11:48:49"
$ws.Range("E104").Value = "Manual"
$ws.Range("A105").Value = "211121"
$ws.Range("D105").Value = "11:49:14"
$ws.Range("E105").Value = "Manual"
$ws.Range("A106").Value = "201252"
$ws.Range("D106").Value = "11:49:29"
$ws.Range("E106").Value = "Manual"
$ws.Range("A107").Value = "201023"
$ws.Range("D107").Value = "11:49:42"
$ws.Range("E107").Value = "Manual"
$ws.Range("A108").Value = "201253"
$ws.Range("D108").Value = "11:49:57"
$ws.Range("E108").Value = "Manual"
$ws.Range("A109").Value = "201337"
$ws.Range("D109").Value = "11:50:11"
$ws.Range("E109").Value = "Manual"
$ws.Range("A110").Value = "201572"
$ws.Range("D110").Value = "11:50:20"
$ws.Range("E110").Value = "Scan"
$ws.Range("A111").Value = "201398"
$ws.Range("D111").Value = "11:50:41"
$ws.Range("E111").Value = "Manual"

# Remove the now-duplicate last row (112); its data already moved into
# row 111 above, so the sheet shrinks from A1:F112 to A1:F111.
$ws.Rows.Item(112).Delete()
